$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for B2:I35 (one PowerShell array literal per spreadsheet row)
$rows = @(
    @(0,0,0,0,0,0,0.008665511265164644,0.01290201944652205),
    @(0,0,0.0627340823970038,0,0,0,0.024263431542461,0.1159311892296172),
    @(0.04270613107822414,0,0,0,0,0,0.0121317157712305,0),
    @(0,0.1647727272727272,0,0,0,0,0.6169844020797205,0),
    @(0,0,0,0,0,0,0.01039861351819757,0.01925953627524309),
    @(0.02832980972515859,0,0,0,0,0.05063291139240502,0,0),
    @(0,0,0.4784644194756565,0,0.74074074074074,0,0,0.02935676888556477),
    @(0.04608879492600428,0,0,0,0,0.1511075949367092,0,0),
    @(0.01564482029598309,0,0,0.2564102564102564,0,0.08306962025316458,0,0),
    @(0.001268498942917547,0,0,0,0,0.01819620253164558,0,0),
    @(0.003805496828752641,0,0,0.05128205128205128,0,0.04430379746835441,0,0),
    @(0,0,0,0,0,0,0.02599653379549393,0.008414360508601349),
    @(0,0.4999999999999997,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0),
    @(0.1238900634249463,0,0,0.05128205128205128,0,0.03243670886075951,0.03466204506065857,0.07068062827225134),
    @(0,0,0,0,0,0,0,0.01869857890800299),
    @(0,0,0,0.02564102564102564,0,0.0134493670886076,0,0),
    @(0,0.005681818181818182,0,0,0,0,0.121317157712305,0),
    @(0.04735729386892183,0,0,0,0,0.009493670886075951,0.008665511265164644,0.06207928197457021),
    @(0.1196617336152212,0,0,0,0,0.01661392405063292,0.05199306759098785,0.03328347045624542),
    @(0,0,0,0,0,0,0,0),
    @(0,0,0.4588014981273418,0,0,0,0.04679376083188907,0.06451009723261054),
    @(0,0,0,0.1282051282051282,0,0.02294303797468356,0,0),
    @(0.04524312896405924,0,0,0,0,0.01424050632911393,0,0),
    @(0,0,0,0,0,0.0007911392405063291,0,0),
    @(0.0008456659619450317,0,0,0,0,0.0411392405063291,0,0),
    @(0.04735729386892183,0,0,0,0,0.0007911392405063291,0,0),
    @(0.006765327695560251,0,0,0,0,0.0379746835443038,0,0),
    @(0.01818181818181819,0,0,0,0,0.05142405063291135,0,0),
    @(0.01606765327695561,0,0,0,0,0.05933544303797461,0,0),
    @(0.244820295983084,0,0,0,0,0,0.02772963604852685,0.0001869857890800299),
    @(0.01522198731501058,0,0,0.4871794871794873,0,0.2602848101265832,0,0),
    @(0,0,0,0,0,0,0,0.008040388930441291),
    @(0,0,0,0,0,0.04905063291139237,0,0)
)

$nRows = $rows.Count
$nCols = 8
$block = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $block[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("B2:I35").Value = $block

# The five "Joint regime area ..." rows (36-40) were dropped entirely
$ws.Range("A36:I40").Delete()

Write-Output "Applied weather-cluster refresh: B2:I35 updated, rows 36-40 removed."
